$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks numeric/percentage while keeping the
# cell as plain text (matching the source data which stores these as text).
# We briefly force a Text number format so the Value assignment is not
# auto-coerced to a number/percentage, then restore the cell's original
# (default) style by copying it from an untouched cell in the same column
# family so no stray formatting is left behind.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $ws.Range("D28").Style
}

Set-TextValue "D2" "294.91"
Set-TextValue "E2" "1.66%"
Set-TextValue "D3" "31.02"
Set-TextValue "E3" "0.88%"
Set-TextValue "D4" "4.916"
Set-TextValue "E4" "-0.59%"
Set-TextValue "D5" "0.07445"
Set-TextValue "E5" "4.03%"
Set-TextValue "D6" "2.219"
Set-TextValue "E6" "23.30%"
Set-TextValue "D7" "7.758"
Set-TextValue "E7" "1.40%"
Set-TextValue "D8" "3.743"
Set-TextValue "E8" "0.11%"
Set-TextValue "D9" "0.9162"
Set-TextValue "E9" "2.02%"
Set-TextValue "D10" "0.08949"
Set-TextValue "E10" "16.83%"
Set-TextValue "D11" "0.1714"
Set-TextValue "E11" "4.21%"
Set-TextValue "D12" "0.08328"
Set-TextValue "E12" "3.89%"
Set-TextValue "D13" "0.03108"
Set-TextValue "E13" "2.59%"
Set-TextValue "D14" "0.1007"
Set-TextValue "E14" "0.52%"
Set-TextValue "D15" "0.001517"
Set-TextValue "E15" "1.44%"
Set-TextValue "D16" "0.005686"
Set-TextValue "E16" "-3.32%"
Set-TextValue "D17" "3.505"
Set-TextValue "E17" "0.89%"
Set-TextValue "D18" "2.076"
Set-TextValue "E18" "-0.32%"
Set-TextValue "E19" "1.58%"
Set-TextValue "E20" "-1.05%"
Set-TextValue "D21" "3.982"
Set-TextValue "E21" "-1.46%"
Set-TextValue "E22" "5.15%"
Set-TextValue "D23" "0.04562"
Set-TextValue "E23" "1.05%"
Set-TextValue "D24" "0.001213"
Set-TextValue "E24" "-0.01%"
Set-TextValue "D25" "0.004621"
Set-TextValue "E25" "15.48%"
Set-TextValue "D26" "0.0001301"
Set-TextValue "E26" "4.13%"
Set-TextValue "D27" "0.0003394"
Set-TextValue "E27" "-95.49%"
Set-TextValue "E39" "0.18%"
Set-TextValue "D40" "0.04480"
Set-TextValue "E40" "2.81%"
Set-TextValue "D41" "0.007294"
Set-TextValue "E41" "-0.52%"
Set-TextValue "D42" "0.008947"
Set-TextValue "D43" "0.1328"
Set-TextValue "E43" "1.84%"
Set-TextValue "D44" "0.001964"
Set-TextValue "E44" "-4.64%"
Set-TextValue "D45" "0.008606"
Set-TextValue "E45" "-6.02%"
Set-TextValue "D46" "0.00006075"
Set-TextValue "E46" "1.58%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "0.05%"
Set-TextValue "D48" "2.230"
Set-TextValue "E48" "-0.69%"
Set-TextValue "D49" "0.002000"
Set-TextValue "E49" "-33.30%"
Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "0.05%"
Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "0.05%"
